$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-22 02:20:57"
$wsZhCn.Range("G4").Value = "2016-01-22 02:21:45"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-22 02:21:10"
$wsDeDe.Range("G4").Value = "2016-01-22 02:22:08"
